# Updated sprint 1 schedule
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Sprint 1 block total (row 4) ---
# The "Total" cell for the Sprint 1 header row no longer carries a value.
$ws.Range("G4").ClearContents()

# --- Row 5: 1.1 User Stories ---
$ws.Range("E5").Value = 1
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 3

# --- Row 6: 1.3 ER Diagram ---
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = 2
$ws.Range("G6").Value = 1
$ws.Range("H6").Value = 8

# --- Row 7: 1.3 Cost Estimation ---
$ws.Range("D7").Value = 1
$ws.Range("G7").Value = 2
$ws.Range("H7").Value = 6

# --- Row 8: 1.4 Proposed Screens ---
$ws.Range("C8").Value = 1
$ws.Range("E8").Value = 2
$ws.Range("G8").Value = 2
$ws.Range("H8").Value = 9

# --- Row 9: 1.5 Schedule Management ---
$ws.Range("B9").Value = 1
$ws.Range("G9").Value = 1
$ws.Range("H9").Value = 9

# --- Row 10: Sprint 2 header total moved from column G to column H ---
$ws.Range("G10").ClearContents()
$ws.Range("H10").Value = 35

# --- View state: scroll the window so column B is left-most and select H10 ---
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("H10").Select()
